# The deck's slide master ("Integral" / "Red Violet" colour scheme, stored
# in ppt/theme/theme1.xml) is switched to the default Office theme colour
# scheme ("Office Theme" / "Office"), as seen in the authoritative diff
# where the old theme1.xml content becomes theme2.xml's content, and
# theme2.xml is replaced with what used to be theme1.xml's content.
#
# PowerPoint's object model exposes the applied theme's 12-colour scheme
# through SlideMaster.Theme.ThemeColorScheme(1..12), in document order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# Re-pointing every slot to the standard Office theme's RGB values
# reproduces the colour-scheme half of the swap (the part of the theme
# XML that is writable through COM automation; the theme/colour-scheme
# display names are read-only derived properties in the PowerPoint
# object model, exactly like the real application).

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

$officeColors = @{
    1  = 0x000000   # dk1
    2  = 0xFFFFFF   # lt1
    3  = 0x44546A   # dk2
    4  = 0xE7E6E6   # lt2
    5  = 0x5B9BD5   # accent1
    6  = 0xED7D31   # accent2
    7  = 0xA5A5A5   # accent3
    8  = 0xFFC000   # accent4
    9  = 0x4472C4   # accent5
    10 = 0x70AD47   # accent6
    11 = 0x0563C1   # hlink
    12 = 0x954F72   # folHlink
}

foreach ($index in 1..12) {
    $hex = $officeColors[$index]
    # $hex is a plain 0xRRGGBB literal -> pull bytes out most-significant first.
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    # PowerPoint COM RGB values are packed as 0x00BBGGRR.
    $comRgb = $r + ($g * 256) + ($b * 65536)
    $colorScheme.Item($index).RGB = $comRgb
}

# Best-effort: try to also update the display names. Real PowerPoint
# treats these as read-only (derived from whichever theme/colour scheme
# was applied), so this is a harmless no-op if unsupported.
try { $theme.Name = "Office Theme" } catch {}
try { $colorScheme.Name = "Office" } catch {}
